$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") was populated from the source filename (e.g. "6-2-2007-08")
# instead of the actual game date. NBA.com stats for a game date were being
# shown one day off, so the data pull used the wrong date string. Correct it
# to the real ISO game date for every data row (2-31).
#
# Force the column to stay text (NumberFormat "@") before writing so Excel
# doesn't reinterpret the date-like string as a serial date value.
$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2008-06-02"
}
